# Apply updated "dSF" (column F) values as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4  = -2
    6  = 2
    19 = 5
    22 = -7
    25 = -3
    29 = -1
    36 = -3
    43 = 3
    45 = 0
    46 = 3
    48 = -3
    49 = -2
    52 = -1
    53 = -1
    55 = -1
    57 = 3
    61 = -7
    66 = -6
    71 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
